# Incorporated reading excel data to create book
# Adds a new "AddBook" test case row (with an extra data column) to the
# "main" worksheet, mirroring a data-driven test case that exercises
# creating a book via ISBN number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")
$ws.Activate()

# Duplicate the "data3" header into the new 5th column (E1).
$ws.Range("E1").Value = "data3"

# New test case row: AddBook / Excel Driven Book / isbnNum / 98761 / Maruf RA
$ws.Range("A6").Value = "AddBook"
$ws.Range("B6").Value = "Excel Driven Book"
$ws.Range("C6").Value = "isbnNum"
$ws.Range("D6").Value = 98761
$ws.Range("E6").Value = "Maruf RA"

# Leave the selection on E2, matching the saved view state.
$ws.Range("E2").Select()
